$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

for ($r = 1; $r -le 18; $r++) {
    $ws.Rows.Item($r).RowHeight = 19.5
}
